# Auto-generated edit script: updates the cryptos price/volume table
# to reflect the latest scrape (reordered rows + refreshed numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '68.817.70'
$ws.Cells.Item(2, 5).Value = '  +1.14%  '
$ws.Cells.Item(3, 4).Value = '3.294.48'
$ws.Cells.Item(3, 5).Value = '  +1.38%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '587.75'
$ws.Cells.Item(5, 5).Value = '  +0.94%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '187.98'
$ws.Cells.Item(6, 5).Value = '  +1.62%  '
$ws.Cells.Item(7, 5).Value = '  +0.02%  '
$ws.Cells.Item(8, 5).Value = '  +1.04%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.131'
$ws.Cells.Item(9, 5).Value = '  -0.47%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '6.67'
$ws.Cells.Item(10, 5).Value = '  +0.45%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.413'
$ws.Cells.Item(11, 5).Value = '  -1.43%  '
$ws.Cells.Item(12, 4).Value = '3.864.76'
$ws.Cells.Item(12, 5).Value = '  +1.40%  '
$ws.Cells.Item(13, 5).Value = '  +1.26%  '
$ws.Cells.Item(14, 2).Value = 'Avalanche'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '27.70'
$ws.Cells.Item(14, 5).Value = '  -1.15%  '
$ws.Cells.Item(15, 2).Value = 'WrappedBTC'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(15, 4).Value = '68.790.22'
$ws.Cells.Item(15, 5).Value = '  +1.09%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.0000170'
$ws.Cells.Item(16, 5).Value = '  -0.20%  '
$ws.Cells.Item(17, 4).Value = '3.286.69'
$ws.Cells.Item(17, 5).Value = '  +2.09%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '5.77'
$ws.Cells.Item(18, 5).Value = '  -0.89%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '13.55'
$ws.Cells.Item(19, 5).Value = '  +0.02%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '420.10'
$ws.Cells.Item(20, 5).Value = '  +6.60%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '7.63'
$ws.Cells.Item(21, 5).Value = '  -0.34%  '
$ws.Cells.Item(22, 2).Value = 'Litecoin'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '71.82'
$ws.Cells.Item(22, 5).Value = '  +0.78%  '
$ws.Cells.Item(23, 2).Value = 'Dai'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.999'
$ws.Cells.Item(23, 5).Value = '  -0.18%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '0.514'
$ws.Cells.Item(24, 5).Value = '  -0.70%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.0000119'
$ws.Cells.Item(25, 5).Value = '  +0.09%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.190'
$ws.Cells.Item(26, 5).Value = '  +0.64%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '9.51'
$ws.Cells.Item(27, 5).Value = '  -2.78%  '
$ws.Cells.Item(28, 5).Value = '  +0.48%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.96'
$ws.Cells.Item(29, 5).Value = '  -0.36%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '22.90'
$ws.Cells.Item(30, 5).Value = '  +0.31%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '5.53'
$ws.Cells.Item(31, 5).Value = '  -2.34%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.26'
$ws.Cells.Item(32, 5).Value = '  -0.58%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '6.93'
$ws.Cells.Item(33, 5).Value = '  -2.65%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '164.64'
$ws.Cells.Item(34, 5).Value = '  +1.65%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.47'
$ws.Cells.Item(35, 5).Value = '  -1.51%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.93'
$ws.Cells.Item(36, 5).Value = '  -0.75%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '26.75'
$ws.Cells.Item(37, 5).Value = '  +0.23%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.798'
$ws.Cells.Item(38, 5).Value = '  -2.35%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '4.52'
$ws.Cells.Item(39, 5).Value = '  -1.22%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '6.39'
$ws.Cells.Item(40, 5).Value = '  -1.73%  '
$ws.Cells.Item(41, 4).Value = '2.668.42'
$ws.Cells.Item(41, 5).Value = '  +1.75%  '
$ws.Cells.Item(42, 2).Value = 'dogwifhat'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.47'
$ws.Cells.Item(42, 5).Value = '  -0.12%  '
$ws.Cells.Item(43, 2).Value = 'Hedera'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.0684'
$ws.Cells.Item(43, 5).Value = '  -0.56%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '341.72'
$ws.Cells.Item(44, 5).Value = '  +1.07%  '
$ws.Cells.Item(45, 2).Value = 'OKB'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '40.65'
$ws.Cells.Item(45, 5).Value = '  -0.88%  '
$ws.Cells.Item(46, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '24.84'
$ws.Cells.Item(46, 5).Value = '  -1.66%  '
$ws.Cells.Item(47, 2).Value = 'VeChain'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.0278'
$ws.Cells.Item(47, 5).Value = '  -1.17%  '
$ws.Cells.Item(48, 2).Value = 'ONDO'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.999'
$ws.Cells.Item(48, 5).Value = '  +1.84%  '
$ws.Cells.Item(49, 2).Value = 'Arweave'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '31.58'
$ws.Cells.Item(49, 5).Value = '  +1.04%  '
$ws.Cells.Item(50, 2).Value = 'Cosmos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '6.27'
$ws.Cells.Item(50, 5).Value = '  -1.32%  '
$ws.Cells.Item(51, 2).Value = 'Stellar'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.101'
$ws.Cells.Item(51, 5).Value = '  -0.59%  '
